{"js": "// Replace the 25 two-digit-by-two-digit multiplication problems in the\n// document's table with their new values. Each old value is unique across\n// the document, so a straight text search-and-replace (format preserving)\n// is unambiguous.\nconst replacements = [\n  [\"92\u00d722=\", \"48\u00d792=\"],\n  [\"44\u00d718=\", \"96\u00d741=\"],\n  [\"80\u00d773=\", \"50\u00d780=\"],\n  [\"42\u00d712=\", \"93\u00d730=\"],\n  [\"55\u00d775=\", \"34\u00d714=\"],\n  [\"77\u00d766=\", \"57\u00d745=\"],\n  [\"66\u00d722=\", \"22\u00d772=\"],\n  [\"83\u00d753=\", \"19\u00d766=\"],\n  [\"55\u00d764=\", \"70\u00d742=\"],\n  [\"96\u00d721=\", \"73\u00d726=\"],\n  [\"43\u00d749=\", \"92\u00d734=\"],\n  [\"37\u00d767=\", \"69\u00d722=\"],\n  [\"92\u00d711=\", \"32\u00d764=\"],\n  [\"46\u00d784=\", \"89\u00d731=\"],\n  [\"33\u00d753=\", \"58\u00d793=\"],\n  [\"35\u00d791=\", \"52\u00d792=\"],\n  [\"30\u00d781=\", \"50\u00d733=\"],\n  [\"86\u00d772=\", \"95\u00d718=\"],\n  [\"96\u00d758=\", \"54\u00d763=\"],\n  [\"78\u00d781=\", \"18\u00d756=\"],\n  [\"29\u00d797=\", \"37\u00d743=\"],\n  [\"99\u00d793=\", \"99\u00d724=\"],\n  [\"91\u00d794=\", \"28\u00d797=\"],\n  [\"30\u00d797=\", \"18\u00d784=\"],\n  [\"84\u00d719=\", \"50\u00d783=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 two-digit-by-two-digit multiplication problems in the\n# document's table with their new values. Each old value is unique across\n# the document, so a straight Find/Replace (format preserving, since\n# Find.Execute only swaps the matched text run) is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{old = \"92\u00d722=\"; new = \"48\u00d792=\"},\n    @{old = \"44\u00d718=\"; new = \"96\u00d741=\"},\n    @{old = \"80\u00d773=\"; new = \"50\u00d780=\"},\n    @{old = \"42\u00d712=\"; new = \"93\u00d730=\"},\n    @{old = \"55\u00d775=\"; new = \"34\u00d714=\"},\n    @{old = \"77\u00d766=\"; new = \"57\u00d745=\"},\n    @{old = \"66\u00d722=\"; new = \"22\u00d772=\"},\n    @{old = \"83\u00d753=\"; new = \"19\u00d766=\"},\n    @{old = \"55\u00d764=\"; new = \"70\u00d742=\"},\n    @{old = \"96\u00d721=\"; new = \"73\u00d726=\"},\n    @{old = \"43\u00d749=\"; new = \"92\u00d734=\"},\n    @{old = \"37\u00d767=\"; new = \"69\u00d722=\"},\n    @{old = \"92\u00d711=\"; new = \"32\u00d764=\"},\n    @{old = \"46\u00d784=\"; new = \"89\u00d731=\"},\n    @{old = \"33\u00d753=\"; new = \"58\u00d793=\"},\n    @{old = \"35\u00d791=\"; new = \"52\u00d792=\"},\n    @{old = \"30\u00d781=\"; new = \"50\u00d733=\"},\n    @{old = \"86\u00d772=\"; new = \"95\u00d718=\"},\n    @{old = \"96\u00d758=\"; new = \"54\u00d763=\"},\n    @{old = \"78\u00d781=\"; new = \"18\u00d756=\"},\n    @{old = \"29\u00d797=\"; new = \"37\u00d743=\"},\n    @{old = \"99\u00d793=\"; new = \"99\u00d724=\"},\n    @{old = \"91\u00d794=\"; new = \"28\u00d797=\"},\n    @{old = \"30\u00d797=\"; new = \"18\u00d784=\"},\n    @{old = \"84\u00d719=\"; new = \"50\u00d783=\"}\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair.old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair.new\n    $found = $find.Execute($pair.old, $false, $true, $false, $false, $false, $true, 1, $false, $pair.new, 2)\n    if (-not $found) {\n        throw \"Text not found: $($pair.old)\"\n    }\n}\n"}
